# "recuperando datos de planilla excel" -- restore/refresh the sample data
# on Hoja1: rename headers, recompute the sensor readings and swap the
# IF() formulas in column D to the new red/yellow/green/blue tagging
# scheme used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "tag"
$ws.Range("D1").Value = "formula"
$ws.Range("F1").Value = "attr"

# --- Row 2: PM_IPA_FERMENTACION_PRESION / stroke ---------------------
$ws.Range("B2").Value = 3.15
$ws.Range("D2").Formula = '=IF(B2>C2,"green","blue")'

# --- Row 3: PM_IPA_FERMENTACION_PRESION / fill ------------------------
$ws.Range("B3").Value = 3.15
$ws.Range("D3").Formula = '=IF(B3>C3,"red","yellow")'

# --- Row 4: PM_IPA_FERMENTACION_PRESION / text ------------------------
$ws.Range("B4").Value = 3.15
$ws.Range("D4").Formula = '=IF(B4>C4,"red","yellow")'

# --- Row 5: PM_IPA_CENTRIFUGADO_MARCHA / fill -------------------------
# values used to be the text "true"/"true"; now numeric 1/1
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Formula = '=IF(B5=C5,"green","blue")'

# --- Cosmetic: selection moved off the table, tidy zoom ---------------
$ws.Range("G5").Select()
$excel.ActiveWindow.Zoom = 100
